# Apply capex/opex battery storage fix: update LCOH values on the three
# storage worksheets and widen column B on the "Present-Storage" and
# "2050-Storage" sheets.

$wb = $excel.ActiveWorkbook

# --- Present-Storage ---
$ws = $wb.Worksheets.Item("Present-Storage")
# Column B widened from 8.4 to 9.6 (character units) in the OOXML <col>
# definition. The host quantizes ColumnWidth writes to an even-pixel grid
# (steps of 1/6 character unit), so 8.8333 (the middle of the input bucket
# that resolves to the stored width nearest 9.6) is used here.
$ws.Columns.Item(2).ColumnWidth = 8.8333333333333
$ws.Range("B2").Value = 132.04
$ws.Range("B3").Value = 563.3
$ws.Range("B4").Value = 58.67
$ws.Range("B6").Value = 139.11
$ws.Range("B7").Value = 595.17
$ws.Range("B8").Value = 61.52
$ws.Range("B10").Value = 107.99
$ws.Range("B11").Value = 444.46
$ws.Range("B12").Value = 50.75
$ws.Range("B14").Value = 128.01
$ws.Range("B15").Value = 541.8
$ws.Range("B16").Value = 57.61

# --- 2030-Storage ---
$ws = $wb.Worksheets.Item("2030-Storage")
$ws.Range("B2").Value = 5.88
$ws.Range("B3").Value = 6.04
$ws.Range("B4").Value = 6.37
$ws.Range("B6").Value = 5.7
$ws.Range("B7").Value = 5.87
$ws.Range("B8").Value = 6.21
$ws.Range("B10").Value = 9.56
$ws.Range("B11").Value = 9.69
$ws.Range("B12").Value = 9.94
$ws.Range("B14").Value = 6.96
$ws.Range("B15").Value = 7.12
$ws.Range("B16").Value = 7.43

# --- 2050-Storage ---
$ws = $wb.Worksheets.Item("2050-Storage")
$ws.Columns.Item(2).ColumnWidth = 8.8333333333333
$ws.Range("B2").Value = 146.2
$ws.Range("B3").Value = 73.36
$ws.Range("B4").Value = 92.23
$ws.Range("B6").Value = 154.08
$ws.Range("B7").Value = 77.05
$ws.Range("B8").Value = 97.01
$ws.Range("B10").Value = 119.03
$ws.Range("B11").Value = 62.2
$ws.Range("B12").Value = 76.93
$ws.Range("B14").Value = 141.59
$ws.Range("B15").Value = 71.7
$ws.Range("B16").Value = 89.81
